# Update specification tables and examples
#
# Row 7 (BOR, BNOT, BRS, BLS) is removed from the keyword table, and those
# four keywords are relocated into the cells that previously held keywords
# which are being dropped entirely from the spec (GET, RETURNS, INF, CEIL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete last row of the table.
$ws.Rows(7).Delete()

# Relocate the remaining keywords to replace the dropped ones.
$ws.Range("B1").Value = "BNOT"
$ws.Range("E1").Value = "BRS"
$ws.Range("E2").Value = "BLS"
$ws.Range("D6").Value = "BOR"

# Update the selected cell to match the saved view.
[void]$ws.Range("F6").Select()
